$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.517.93'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.952.33'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.615'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.27'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +5.44%  '
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.76%  '
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -7.73%  '
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.26'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +6.56%  '
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.238.22'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.54%  '
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.30%  '
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.50'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.41%  '
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.86%  '
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.948.36'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.46%  '
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.466.20'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.31'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.40%  '
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '228.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.64%  '
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.13%  '
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.45'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.26%  '
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.83%  '
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.143'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +7.64%  '
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.23%  '
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '159.73'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.49%  '
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.28'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.78%  '
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +20.44%  '
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.14%  '
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.72'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.47%  '
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.45'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +7.09%  '
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +10.49%  '
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.28'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +5.41%  '
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'BinanceUSD'
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.77'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.41%  '
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -11.74%  '
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0963'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.00%  '
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.67%  '
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.27%  '
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.62%  '
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.72'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.357.53'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.62%  '
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.50'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.06%  '
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.17%  '
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.73%  '
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.67%  '
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '45.22'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.77%  '
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.134.44'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.77%  '
$ws.Range("E51").Style = "Normal"
